$p = $ppt.ActivePresentation

# --- Slide 2 ("Co je CryptoPeek?"): fix "rychlí" -> "rychlý" -----------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "Aplikace byla vyvinuta z důvodu mé vlastní potřeby, abych měl rychlý a jednoduchý přístup ke zjištění kurzu a novinek"

# --- Slide 3 ("Funkce"): fix "Rychlí přehled" -> "Rychlý přehled" -----------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange

$para3a = $tr3.Paragraphs(1, 1)
$run3a = $para3a.Runs(1, 1)
$run3a.Text = "Rychlý přehled "

# --- Slide 3: fix "které byli uvedeny" -> "které byly uvedeny" -------------
# This text lives in a run that is grouped together (by COM, via identical
# formatting) with the following "coinmarketcap" / ")" runs, so we target the
# exact character span of the original run (not the whole merged COM run) to
# replace it in place without disturbing the sibling runs. The span is
# computed (rather than hard-coded) from the paragraph's known original text.
$para3b = $tr3.Paragraphs(2, 1)
$para2PlainText = "Vyhledávání aktuálního kurzu zadané kryptoměny (vyhledá všechny měny, které byli uvedeny na coinmarketcap)"
$oldRunText = "(vyhledá všechny měny, které byli uvedeny na "
$newRunText = "(vyhledá všechny měny, které byly uvedeny na "
$startPos = $para2PlainText.IndexOf($oldRunText) + 1
$chars3b = $para3b.Characters($startPos, $oldRunText.Length)
$chars3b.Text = $newRunText
